# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-01-06 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-01-07 Tuesday", 2)

# Update the division problems in the table, cell by cell, to avoid any
# ambiguity from duplicate / cross-colliding values when using Find/Replace.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "15÷9=" },
    @{ Row = 1;  Col = 2; Text = "34÷7=" },
    @{ Row = 1;  Col = 3; Text = "41÷6=" },
    @{ Row = 1;  Col = 4; Text = "65÷6=" },
    @{ Row = 1;  Col = 5; Text = "17÷8=" },

    @{ Row = 5;  Col = 1; Text = "16÷2=" },
    @{ Row = 5;  Col = 2; Text = "52÷6=" },
    @{ Row = 5;  Col = 3; Text = "52÷6=" },
    @{ Row = 5;  Col = 4; Text = "63÷4=" },
    @{ Row = 5;  Col = 5; Text = "64÷8=" },

    @{ Row = 9;  Col = 1; Text = "44÷5=" },
    @{ Row = 9;  Col = 2; Text = "12÷2=" },
    @{ Row = 9;  Col = 3; Text = "59÷7=" },
    @{ Row = 9;  Col = 4; Text = "94÷6=" },
    @{ Row = 9;  Col = 5; Text = "70÷6=" },

    @{ Row = 13; Col = 1; Text = "23÷3=" },
    @{ Row = 13; Col = 2; Text = "46÷5=" },
    @{ Row = 13; Col = 3; Text = "80÷6=" },
    @{ Row = 13; Col = 4; Text = "29÷5=" },
    @{ Row = 13; Col = 5; Text = "46÷6=" },

    @{ Row = 17; Col = 1; Text = "44÷3=" },
    @{ Row = 17; Col = 2; Text = "33÷3=" },
    @{ Row = 17; Col = 3; Text = "64÷7=" },
    @{ Row = 17; Col = 4; Text = "96÷6=" },
    @{ Row = 17; Col = 5; Text = "58÷2=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
